$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.448
$ws.Range("G2").Value = -6.9
$ws.Range("H2").Value = -6.9
$ws.Range("I2").Value = -11.7
$ws.Range("J2").Value = -11.7
$ws.Range("K2").Value = -0.239
$ws.Range("L2").Value = -11.95
$ws.Range("M2").Value = 0.221
$ws.Range("N2").Value = 0.04446680080482898
$ws.Range("O2").Value = -0.9246861924686193
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0.221
$ws.Range("U2").Value = 0.02
$ws.Range("V2").Value = 0.004024144869215292
$ws.Range("W2").Value = -0.5431818181818182
$ws.Range("X2").Value = 0.05389978651585595
$ws.Range("Y2").Value = -0.5970816046976741
$ws.Range("Z2").Value = 0.05012531328320802
$ws.Range("AA2").Value = -0.5864661654135338
$ws.Range("AB2").Value = 0.05389978651585595
$ws.Range("AC2").Value = -0.6403659519293897
$ws.Range("AG2").Value = -0.02
$ws.Range("AJ2").Value = -0.00404040404040404
$ws.Range("AK2").Value = -0.03610108303249098
$ws.Range("AL2").Value = 0.008
$ws.Range("AM2").Value = 0.008
$ws.Range("AN2").Value = -0
$ws.Range("AO2").Value = -29.25
$ws.Range("AP2").Value = 0.09009009009009009
$ws.Range("AQ2").Value = -29.25

# Row 3
$ws.Range("D3").Value = -0.448
$ws.Range("G3").Value = -6.9
$ws.Range("H3").Value = -6.9
$ws.Range("I3").Value = -11.7
$ws.Range("J3").Value = -11.7
$ws.Range("K3").Value = -0.239
$ws.Range("L3").Value = -11.95
$ws.Range("M3").Value = 0.221
$ws.Range("N3").Value = 0.04446680080482898
$ws.Range("O3").Value = -0.9246861924686193
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.221
$ws.Range("U3").Value = 0.02
$ws.Range("V3").Value = 0.004024144869215292
$ws.Range("W3").Value = -0.5431818181818182
$ws.Range("X3").Value = 0.05389978651585595
$ws.Range("Y3").Value = -0.5970816046976741
$ws.Range("Z3").Value = 0.05012531328320802
$ws.Range("AA3").Value = -0.5864661654135338
$ws.Range("AB3").Value = 0.05389978651585595
$ws.Range("AC3").Value = -0.6403659519293897
$ws.Range("AG3").Value = -0.02
$ws.Range("AJ3").Value = -0.00404040404040404
$ws.Range("AK3").Value = -0.03610108303249098
$ws.Range("AL3").Value = 0.008
$ws.Range("AM3").Value = 0.008
$ws.Range("AN3").Value = -0
$ws.Range("AO3").Value = -29.25
$ws.Range("AP3").Value = 0.09009009009009009
$ws.Range("AQ3").Value = -29.25
